$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "614×7=4298" "282×3=846"
Replace-Text "721×9=6489" "296×7=2072"
Replace-Text "731×2=1462" "569×7=3983"
Replace-Text "384×7=2688" "867×9=7803"
Replace-Text "332×3=996" "378×2=756"
Replace-Text "790×8=6320" "478×2=956"
Replace-Text "280×4=1120" "923×2=1846"
Replace-Text "384×2=768" "825×3=2475"
Replace-Text "370×3=1110" "710×2=1420"
Replace-Text "761×8=6088" "145×8=1160"
Replace-Text "743×8=5944" "525×8=4200"
Replace-Text "208×8=1664" "954×9=8586"
Replace-Text "140×4=560" "710×9=6390"
Replace-Text "179×3=537" "905×6=5430"
Replace-Text "812×9=7308" "434×6=2604"
Replace-Text "166×2=332" "224×2=448"
Replace-Text "614×3=1842" "936×4=3744"
Replace-Text "330×5=1650" "784×6=4704"
Replace-Text "838×8=6704" "668×7=4676"
Replace-Text "374×7=2618" "399×2=798"
Replace-Text "266×9=2394" "900×4=3600"
Replace-Text "804×6=4824" "640×2=1280"
Replace-Text "545×7=3815" "356×9=3204"
Replace-Text "397×5=1985" "250×5=1250"
Replace-Text "973×3=2919" "563×2=1126"
